$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (index 1)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 6).Value = 523
$ws1.Cells.Item(3, 6).Value = 105
$ws1.Cells.Item(8, 6).Value = 1105
$ws1.Cells.Item(9, 6).Value = 226
$ws1.Cells.Item(10, 6).Value = 171
$ws1.Cells.Item(11, 6).Value = 277
$ws1.Cells.Item(12, 6).Value = 1737
$ws1.Cells.Item(13, 6).Value = 643
$ws1.Cells.Item(14, 6).Value = 298
$ws1.Cells.Item(15, 6).Value = 379
$ws1.Cells.Item(16, 6).Value = 3916
$ws1.Cells.Item(18, 6).Value = 433
$ws1.Cells.Item(21, 6).Value = 1247
$ws1.Cells.Item(23, 6).Value = 1894
$ws1.Cells.Item(24, 6).Value = 2910
$ws1.Cells.Item(25, 6).Value = 1733
$ws1.Cells.Item(26, 6).Value = 96
$ws1.Cells.Item(27, 6).Value = 33
$ws1.Cells.Item(28, 6).Value = 141
$ws1.Cells.Item(32, 6).Value = 2120
$ws1.Cells.Item(33, 6).Value = 950
$ws1.Cells.Item(34, 6).Value = 2198
$ws1.Cells.Item(36, 6).Value = 543
$ws1.Cells.Item(37, 6).Value = 363
$ws1.Cells.Item(39, 6).Value = 862
$ws1.Cells.Item(40, 6).Value = 373
$ws1.Cells.Item(41, 6).Value = 1016
$ws1.Cells.Item(42, 6).Value = 849
$ws1.Cells.Item(43, 6).Value = 1116
$ws1.Cells.Item(44, 6).Value = 256
$ws1.Cells.Item(45, 6).Value = 479
$ws1.Cells.Item(46, 6).Value = 326
$ws1.Cells.Item(47, 6).Value = 262
$ws1.Cells.Item(48, 6).Value = 3423

# Sheet 2: 演出 (index 2)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(6, 6).Value = 19
$ws2.Cells.Item(11, 6).Value = 855
$ws2.Cells.Item(20, 6).Value = 23

# Sheet 4: 全部类型 (index 4)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(2, 6).Value = 523
$ws4.Cells.Item(3, 6).Value = 105
$ws4.Cells.Item(7, 6).Value = 1105
$ws4.Cells.Item(8, 6).Value = 226
$ws4.Cells.Item(9, 6).Value = 171
$ws4.Cells.Item(11, 6).Value = 1737
$ws4.Cells.Item(12, 6).Value = 643
$ws4.Cells.Item(13, 6).Value = 298
$ws4.Cells.Item(14, 6).Value = 379
$ws4.Cells.Item(15, 6).Value = 3916
$ws4.Cells.Item(18, 6).Value = 19
$ws4.Cells.Item(21, 6).Value = 1247
$ws4.Cells.Item(22, 6).Value = 2910
$ws4.Cells.Item(24, 6).Value = 1733
$ws4.Cells.Item(25, 6).Value = 96
$ws4.Cells.Item(26, 6).Value = 33
$ws4.Cells.Item(28, 6).Value = 141
$ws4.Cells.Item(29, 6).Value = 855
$ws4.Cells.Item(33, 6).Value = 2120
$ws4.Cells.Item(35, 6).Value = 950
$ws4.Cells.Item(36, 6).Value = 2198
$ws4.Cells.Item(37, 6).Value = 543
$ws4.Cells.Item(38, 6).Value = 363
$ws4.Cells.Item(39, 6).Value = 862
$ws4.Cells.Item(40, 6).Value = 1016
$ws4.Cells.Item(41, 6).Value = 849
$ws4.Cells.Item(42, 6).Value = 1116
$ws4.Cells.Item(43, 6).Value = 256
$ws4.Cells.Item(44, 6).Value = 479
$ws4.Cells.Item(45, 6).Value = 326
$ws4.Cells.Item(46, 6).Value = 23
$ws4.Cells.Item(47, 6).Value = 262
$ws4.Cells.Item(48, 6).Value = 3423
